$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2511.8235
$ws.Range("I40").Value = 2083.1667
$ws.Range("J40").Value = 2745.6365
$ws.Range("K40").Value = 2083.1667
$ws.Range("L40").Value = 2745.6365
$ws.Range("M40").Value = -1908.1667
$ws.Range("N40").Value = -3095.6365

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3925.6191
$ws.Range("I74").Value = 3780.7693
$ws.Range("J74").Value = 4161
$ws.Range("K74").Value = 3780.7693
$ws.Range("L74").Value = 4161
$ws.Range("M74").Value = -2844.7693
$ws.Range("N74").Value = -6033

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3925.6191
$ws.Range("I77").Value = 3780.7693
$ws.Range("J77").Value = 4161
$ws.Range("K77").Value = 18903.8465
$ws.Range("L77").Value = 20805
$ws.Range("M77").Value = -14223.8465
$ws.Range("N77").Value = -30165

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 10181.333
$ws.Range("I86").Value = 12426
$ws.Range("J86").Value = 5692
$ws.Range("K86").Value = 12426
$ws.Range("L86").Value = 5692
$ws.Range("M86").Value = -11303
$ws.Range("N86").Value = -7938

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 46166.668
$ws.Range("J87").Value = 46166.668
$ws.Range("L87").Value = 46166.668
$ws.Range("N87").Value = -48662.668

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 10181.333
$ws.Range("I89").Value = 12426
$ws.Range("J89").Value = 5692
$ws.Range("K89").Value = 62130
$ws.Range("L89").Value = 28460
$ws.Range("M89").Value = -56514
$ws.Range("N89").Value = -39692

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 46166.668
$ws.Range("J90").Value = 46166.668
$ws.Range("L90").Value = 138500.004
$ws.Range("N90").Value = -150980.004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3176.2173
$ws.Range("I2").Value = 1282.2142
$ws.Range("J2").Value = 6122.4443
$ws.Range("K2").Value = 1282.2142
$ws.Range("L2").Value = 6122.4443
$ws.Range("M2").Value = -1169.2142
$ws.Range("N2").Value = -6348.4443

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4142.0625
$ws.Range("I45").Value = 3262.3
$ws.Range("J45").Value = 5608.3335
$ws.Range("K45").Value = 3262.3
$ws.Range("L45").Value = 5608.3335
$ws.Range("M45").Value = -2885.3
$ws.Range("N45").Value = -6362.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2524.125
$ws.Range("I61").Value = 2025.7333
$ws.Range("K61").Value = 2025.7333
$ws.Range("M61").Value = -1813.7333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 26986.75
$ws.Range("J80").Value = 26986.75
$ws.Range("L80").Value = 26986.75
$ws.Range("N80").Value = -28982.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 26986.75
$ws.Range("J83").Value = 26986.75
$ws.Range("L83").Value = 80960.25
$ws.Range("N83").Value = -90944.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 3960
$ws.Range("I97").Value = 3000
$ws.Range("J97").Value = 4600
$ws.Range("K97").Value = 3000
$ws.Range("L97").Value = 4600
$ws.Range("M97").Value = -2504
$ws.Range("N97").Value = -5592

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 200002600
$ws.Range("I102").Value = 3250
$ws.Range("K102").Value = 3250
$ws.Range("M102").Value = -1628

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3176.2173
$ws.Range("I116").Value = 1282.2142
$ws.Range("J116").Value = 6122.4443
$ws.Range("K116").Value = 1282.2142
$ws.Range("L116").Value = 6122.4443
$ws.Range("M116").Value = 1011.7858
$ws.Range("N116").Value = -10710.4443

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2524.125
$ws.Range("I136").Value = 2025.7333
$ws.Range("K136").Value = 6077.199900000001
$ws.Range("M136").Value = -3527.199900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3176.2173
$ws.Range("I3").Value = 1282.2142
$ws.Range("J3").Value = 6122.4443
$ws.Range("K3").Value = 1282.2142
$ws.Range("L3").Value = 6122.4443
$ws.Range("M3").Value = -1168.2142
$ws.Range("N3").Value = -6350.4443

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18898.166
$ws.Range("J82").Value = 20626.4
$ws.Range("L82").Value = 20626.4
$ws.Range("N82").Value = -21392.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 18898.166
$ws.Range("J85").Value = 20626.4
$ws.Range("L85").Value = 20626.4
$ws.Range("N85").Value = -23278.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3217.2727
$ws.Range("I99").Value = 1898
$ws.Range("J99").Value = 4316.6665
$ws.Range("K99").Value = 1898
$ws.Range("L99").Value = 4316.6665
$ws.Range("M99").Value = -400
$ws.Range("N99").Value = -7312.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 7796.44
$ws.Range("I107").Value = 7859.591
$ws.Range("J107").Value = 7333.3335
$ws.Range("K107").Value = 7859.591
$ws.Range("L107").Value = 7333.3335
$ws.Range("M107").Value = -5939.591
$ws.Range("N107").Value = -11173.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 47105.883
$ws.Range("J126").Value = 47105.883
$ws.Range("L126").Value = 47105.883
$ws.Range("N126").Value = -56985.883

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 37455.5
$ws.Range("J132").Value = 37455.5
$ws.Range("L132").Value = 37455.5
$ws.Range("N132").Value = -47575.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 8481.223
$ws.Range("I50").Value = 4000
$ws.Range("K50").Value = 4000
$ws.Range("M50").Value = -3375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 7535
$ws.Range("J51").Value = 10570
$ws.Range("L51").Value = 10570
$ws.Range("N51").Value = -12042

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 7535
$ws.Range("J61").Value = 10570
$ws.Range("L61").Value = 10570
$ws.Range("N61").Value = -11266

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2853.3333
$ws.Range("I105").Value = 1473.3334
$ws.Range("J105").Value = 5613.3335
$ws.Range("K105").Value = 1473.3334
$ws.Range("L105").Value = 5613.3335
$ws.Range("M105").Value = 273.6666
$ws.Range("N105").Value = -9107.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1593.3334
$ws.Range("I107").Value = 1590
$ws.Range("J107").Value = 1600
$ws.Range("K107").Value = 1590
$ws.Range("L107").Value = 1600
$ws.Range("M107").Value = 330
$ws.Range("N107").Value = -5440

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 20634
$ws.Range("J109").Value = 20634
$ws.Range("L109").Value = 20634
$ws.Range("N109").Value = -22714

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 2973.8438
$ws.Range("J117").Value = 3540.8845
$ws.Range("L117").Value = 10622.6535
$ws.Range("N117").Value = -17506.6535

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 25645142
$ws.Range("I118").Value = 66667388
$ws.Range("J118").Value = 6238
$ws.Range("K118").Value = 200002164
$ws.Range("L118").Value = 18714
$ws.Range("M118").Value = -200000921
$ws.Range("N118").Value = -21200

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2841.34
$ws.Range("I122").Value = 210.05556
$ws.Range("J122").Value = 3418.939
$ws.Range("K122").Value = 1890.50004
$ws.Range("L122").Value = 30770.451
$ws.Range("M122").Value = 559.4999599999999
$ws.Range("N122").Value = -35670.451

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1698.125
$ws.Range("J129").Value = 2665.7144
$ws.Range("L129").Value = 7997.1432
$ws.Range("N129").Value = -17997.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 16980.727
$ws.Range("J57").Value = 16980.727
$ws.Range("L57").Value = 16980.727
$ws.Range("N57").Value = -18620.727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2005.7273
$ws.Range("I102").Value = 1848.5883
$ws.Range("K102").Value = 1848.5883
$ws.Range("M102").Value = -226.5882999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2468.04
$ws.Range("I132").Value = 2180.8125
$ws.Range("K132").Value = 6542.4375
$ws.Range("M132").Value = -4012.4375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 60198.6
$ws.Range("J133").Value = 60198.6
$ws.Range("L133").Value = 60198.6
$ws.Range("N133").Value = -65258.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 29677
$ws.Range("J109").Value = 29677
$ws.Range("L109").Value = 29677
$ws.Range("N109").Value = -32451

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 57535.832
$ws.Range("J135").Value = 57535.832
$ws.Range("L135").Value = 57535.832
$ws.Range("N135").Value = -67675.83199999999
